# NYPD CompStat weekly report refresh: roll the report forward one week
# (Volume 30 Number 52, week of 12/25/2023-12/31/2023) to
# (Volume 31 Number 1, week of 1/1/2024-1/7/2024), bump the "13/30 Year"
# comparison labels, update the '22 vs 'XX historical labels to '23 vs 'XX,
# and refresh all of the Week-to-Date / 28-Day / Year-to-Date / 2-Year /
# 13-Year / 30-Year crime-count tables with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) - edit the specific substrings
# in place so the rest of each string is left untouched. Edits are done
# right-to-left within each cell so earlier character offsets stay valid.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  52" -> "Volume 31   Number  1"
$ws.Range("A8").Characters(21,2).Text = "1"
$ws.Range("A8").Characters(8,2).Text = "31"

# C9: "Report Covering the Week  12/25/2023  Through  12/31/2023"
#  -> "Report Covering the Week  1/1/2024  Through  1/7/2024"
$ws.Range("C9").Characters(48,10).Text = "1/7/2024"
$ws.Range("C9").Characters(27,10).Text = "1/1/2024"

# M12: "13 Year (2010)" -> "14 Year (2010)"
$ws.Range("M12").Characters(1,2).Text = "14"

# N12: "30 Year (1993)" -> "31 Year (1993)"
$ws.Range("N12").Characters(1,2).Text = "31"

# ---------------------------------------------------------------------
# Row 13: Week to Date / 28 Day / Year to Date column-header years
# ---------------------------------------------------------------------
$ws.Range("C13").Value = 2024
$ws.Range("D13").Value = 2023
$ws.Range("F13").Value = 2024
$ws.Range("G13").Value = 2023
$ws.Range("I13").Value = 2024
$ws.Range("J13").Value = 2023

# ---------------------------------------------------------------------
# Rows 14-29: Crime Complaints weekly table values
# ---------------------------------------------------------------------

# Murder
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -66.666666666666
$ws.Range("N14").Value = -87.5

# Rape
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = -60
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = -53.571428571428
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -60
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -60

# Robbery
$ws.Range("C16").Value = 42
$ws.Range("D16").Value = 46
$ws.Range("E16").Value = -8.695652173913
$ws.Range("F16").Value = 197
$ws.Range("G16").Value = 147
$ws.Range("H16").Value = 34.013605442176
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -8.695652173913
$ws.Range("L16").Value = -2.325581395348
$ws.Range("M16").Value = -33.333333333333
$ws.Range("N16").Value = -87.647058823529

# Fel. Assault
$ws.Range("C17").Value = 75
$ws.Range("D17").Value = 84
$ws.Range("E17").Value = -10.714285714285
$ws.Range("F17").Value = 303
$ws.Range("G17").Value = 278
$ws.Range("H17").Value = 8.992805755395
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 84
$ws.Range("K17").Value = -10.714285714285
$ws.Range("L17").Value = 4.166666666666
$ws.Range("M17").Value = 8.695652173913
$ws.Range("N17").Value = -59.459459459459

# Burglary
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 55
$ws.Range("E18").Value = -45.454545454545
$ws.Range("F18").Value = 133
$ws.Range("G18").Value = 177
$ws.Range("H18").Value = -24.858757062146
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 55
$ws.Range("K18").Value = -45.454545454545
$ws.Range("L18").Value = -37.5
$ws.Range("M18").Value = -45.454545454545
$ws.Range("N18").Value = -87.5

# Gr. Larceny
$ws.Range("C19").Value = 82
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = -18
$ws.Range("F19").Value = 387
$ws.Range("G19").Value = 406
$ws.Range("H19").Value = -4.679802955665
$ws.Range("I19").Value = 82
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = -18
$ws.Range("L19").Value = -1.204819277108
$ws.Range("M19").Value = 32.258064516129
$ws.Range("N19").Value = -32.231404958677

# G.L.A.
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = 15.625
$ws.Range("F20").Value = 134
$ws.Range("G20").Value = 130
$ws.Range("H20").Value = 3.076923076923
$ws.Range("I20").Value = 37
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = 15.625
$ws.Range("L20").Value = 2.777777777777
$ws.Range("M20").Value = 48
$ws.Range("N20").Value = -80.526315789473

# TOTAL
$ws.Range("C21").Value = 271
$ws.Range("D21").Value = 328
$ws.Range("E21").Value = -17.378048780487
$ws.Range("F21").Value = 1169
$ws.Range("G21").Value = 1168
$ws.Range("H21").Value = 0.085616438356
$ws.Range("I21").Value = 271
$ws.Range("J21").Value = 328
$ws.Range("K21").Value = -17.378048780487
$ws.Range("L21").Value = -5.244755244755
$ws.Range("M21").Value = -4.240282685512
$ws.Range("N21").Value = -75.228519195612

# Transit
$ws.Range("C22").Value = 5
$ws.Range("E22").Value = -28.571428571428
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 12
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -28.571428571428
$ws.Range("L22").Value = -28.571428571428
$ws.Range("M22").Value = -44.444444444444

# Housing
$ws.Range("D23").Value = 34
$ws.Range("E23").Value = -20.588235294117
$ws.Range("F23").Value = 103
$ws.Range("G23").Value = 112
$ws.Range("H23").Value = -8.035714285714
$ws.Range("I23").Value = 27
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = -20.588235294117
$ws.Range("L23").Value = -10
$ws.Range("M23").Value = 58.823529411764

# Petit Larceny
$ws.Range("C24").Value = 182
$ws.Range("D24").Value = 199
$ws.Range("E24").Value = -8.542713567839
$ws.Range("F24").Value = 917
$ws.Range("G24").Value = 888
$ws.Range("H24").Value = 3.265765765765
$ws.Range("I24").Value = 182
$ws.Range("J24").Value = 199
$ws.Range("K24").Value = -8.542713567839
$ws.Range("L24").Value = 1.111111111111
$ws.Range("M24").Value = 18.181818181818

# Misd. Assault
$ws.Range("C25").Value = 112
$ws.Range("D25").Value = 105
$ws.Range("E25").Value = 6.666666666666
$ws.Range("F25").Value = 450
$ws.Range("G25").Value = 413
$ws.Range("H25").Value = 8.958837772397
$ws.Range("I25").Value = 112
$ws.Range("J25").Value = 105
$ws.Range("K25").Value = 6.666666666666
$ws.Range("L25").Value = 25.842696629213
$ws.Range("M25").Value = -18.840579710144

# UCR Rape*
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -41.666666666666
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = -35.294117647058
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -41.666666666666
$ws.Range("L26").Value = -22.222222222222

# Other Sex Crimes
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 61
$ws.Range("G27").Value = 39
$ws.Range("H27").Value = 56.410256410256
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 55.555555555555

# Shooting Vic.
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 30.769230769230
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -25
$ws.Range("M28").Value = -70
$ws.Range("N28").Value = -93.75

# Shooting Inc.
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 3
$ws.Range("F29").Value = 14
$ws.Range("G29").Value = 13
$ws.Range("H29").Value = 7.692307692307
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -62.5
$ws.Range("N29").Value = -93.023255813953

# ---------------------------------------------------------------------
# Row 30 (Hate Crimes): C30, I30 and J30 go from numeric counts to the
# literal text "0" (shared with D30's existing "0"); K30 goes from a
# numeric % change to the literal text "***.*" (shared with E30/M30/N30).
# Force text entry via a Text number format, then copy the number
# format from a sibling cell that already holds the same text so the
# cell re-uses that cell's style, matching how D30/E30 are styled.
# ---------------------------------------------------------------------
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 600
$ws.Range("L30").Value = -100

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("I30").NumberFormat = "@"
$ws.Range("I30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("I30").PasteSpecial(-4122)

$ws.Range("J30").NumberFormat = "@"
$ws.Range("J30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("J30").PasteSpecial(-4122)

$ws.Range("K30").NumberFormat = "@"
$ws.Range("K30").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("K30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Historical Perspective table: bump the "current year" column from
# 2022 to 2023 and refresh its counts / % changes (rows 36-41, 43).
# Row 42 (G.L.A.) is unchanged. The '22 vs 'XX column-header labels
# are plain shared strings, updated directly to '23 vs 'XX.
# ---------------------------------------------------------------------
$ws.Range("J35").Value = 2023
$ws.Range("K35").Value = "'23 vs '01"
$ws.Range("L35").Value = "'23 vs '98"
$ws.Range("M35").Value = "'23 vs '93"
$ws.Range("N35").Value = "'23 vs '90"

# Murder
$ws.Range("J36").Value = 67
$ws.Range("K36").Value = -59.638554216867
$ws.Range("L36").Value = -50.370370370370
$ws.Range("M36").Value = -85.864978902953
$ws.Range("N36").Value = -86.732673267326

# Rape
$ws.Range("J37").Value = 217
$ws.Range("K37").Value = -40.547945205479
$ws.Range("L37").Value = -52.723311546841
$ws.Range("M37").Value = -63.468013468013
$ws.Range("N37").Value = -69.350282485875

# Robbery
$ws.Range("J38").Value = 2537
$ws.Range("K38").Value = -51.740536427620
$ws.Range("L38").Value = -64.065155807365
$ws.Range("M38").Value = -84.939151083407
$ws.Range("N38").Value = -87.629821054171

# Fel. Assault
$ws.Range("J39").Value = 4232
$ws.Range("K39").Value = -7.375793390238
$ws.Range("L39").Value = -18.942731277533
$ws.Range("M39").Value = -49.558998808104
$ws.Range("N39").Value = -57.586690719583

# Burglary
$ws.Range("J40").Value = 2031
$ws.Range("K40").Value = -52.368667917448
$ws.Range("L40").Value = -62.971741112124
$ws.Range("M40").Value = -83.230121377260
$ws.Range("N40").Value = -86.453678383245

# Gr. Larceny
$ws.Range("J41").Value = 5753
$ws.Range("K41").Value = 46.052297537446
$ws.Range("L41").Value = 55.950121984277
$ws.Range("M41").Value = -16.960161662817
$ws.Range("N41").Value = -34.535730541647

# TOTAL
$ws.Range("J43").Value = 16702
$ws.Range("K43").Value = -23.988531379420
$ws.Range("L43").Value = -35.388781431334
$ws.Range("M43").Value = -69.588492352512
$ws.Range("N43").Value = -75.450149192303
